$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.433.92"
$ws.Range("D3").Value = "1.674.32"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  +0.81%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.06"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5348"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06407"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.04"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07849"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.559"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.671.92"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "1.903.70"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5660"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "0.0₅8197"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.50"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "26.470.76"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.738"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "199.41"
$ws.Range("E21").Value = "  +4.23%  "
$ws.Range("E22").Value = "  +2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.079"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.75"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.263"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.24"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05896"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.287"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.591"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.315"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.619"
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9724"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.852"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5840"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").Value = "1.083.54"
$ws.Range("E40").Value = "  +3.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.934"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8678"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.29"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "1.813.74"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.054"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("E51").Value = "  +0.32%  "
